# Rename the two "open-*" question types to their new, clearer names:
#   open-num  -> numeric
#   open-char -> string
#
# (per commit message: "changed type name from open-num and open-char to
# numeric and string")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Type") holds the values that need renaming.
# Row 3: "How old are you?"                         open-num  -> numeric
# Row 4: "How many years have you been using R?"    open-num  -> numeric
# Row 6: "What is your favourate R package?"         open-char -> string
$ws.Range("C3").Value = "numeric"
$ws.Range("C4").Value = "numeric"
$ws.Range("C6").Value = "string"

# Reflect the author's final selection in the saved sheet view.
$ws.Range("C7").Select()
